$wb = $excel.ActiveWorkbook

# --- Rename sheets (28-01-2024 -> 30-01-2024) ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws1.Name = "inventory 30-01-2024"
$ws2.Name = "transactions 30-01-2024"

# --- Helper: write a text value into a cell without it being coerced to a
#     number (important for barcodes like "045496870775" with a leading
#     zero) and without leaving a stray quote-prefix style behind. ---
function Set-TextCell($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# ================= Sheet 1: inventory =================
# Give column A on the data rows the same bordered/bold style used by the
# header row (B1) - matches s="1" on A2:A6 in the target file.
$ws1.Cells.Item(1, 2).Copy()
$ws1.Range("A2:A6").PasteSpecial(-4122)

$inv = @(
    @(0, "123456",       "Cheetos", 2.5,  32, "Snacks", "Costco", "30-01-2024 14:15:31"),
    @(1, "045496870775", "Doritos", 2,    50, "Snacks", "Costco", "30-01-2024 14:13:23"),
    @(2, "111",          "Test1",   3.45, 39, "Snacks", "Costco", "30-01-2024 14:16:22"),
    @(3, "112",          "Test2",   2.5,  36, "mwe",    "zerg",   "30-01-2024 14:28:21"),
    @(4, "113",          "Test3",   3.5,  44, "mwe",    "zerg",   "30-01-2024 14:28:32")
)

$r = 2
foreach ($row in $inv) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    Set-TextCell $ws1.Cells.Item($r, 2) $row[1]
    Set-TextCell $ws1.Cells.Item($r, 3) $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    Set-TextCell $ws1.Cells.Item($r, 6) $row[5]
    Set-TextCell $ws1.Cells.Item($r, 7) $row[6]
    Set-TextCell $ws1.Cells.Item($r, 8) $row[7]
    $r++
}

# ================= Sheet 2: transactions =================
# Append 3 new transaction rows (7-9); rows 2-6 are unchanged.
$ws2.Cells.Item(1, 2).Copy()
$ws2.Range("A7:A9").PasteSpecial(-4122)

$txn = @(
    @(5, "111", "Test1", 3.45, 3, "Snacks", "Costco", "30-01-2024 14:29:18"),
    @(6, "112", "Test2", 2.5,  4, "mwe",    "zerg",   "30-01-2024 14:29:18"),
    @(7, "113", "Test3", 3.5,  2, "mwe",    "zerg",   "30-01-2024 14:29:18")
)

$r = 7
foreach ($row in $txn) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    Set-TextCell $ws2.Cells.Item($r, 2) $row[1]
    Set-TextCell $ws2.Cells.Item($r, 3) $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    Set-TextCell $ws2.Cells.Item($r, 6) $row[5]
    Set-TextCell $ws2.Cells.Item($r, 7) $row[6]
    Set-TextCell $ws2.Cells.Item($r, 8) $row[7]
    $r++
}

Write-Output "done"
